# Blind model names in the "Evaluations" sheet (column C) so reviewers
# cannot see which underlying model produced each translation.
#
# Mapping (per blinding key):
#   claude-opus-4.5 -> Model A
#   gemini-3-pro    -> Model B
#   gpt-5.1         -> Model C
#   kimi-k2         -> Model D

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Evaluations")

$mapping = @{
    "gpt-5.1"         = "Model C"
    "claude-opus-4.5" = "Model A"
    "gemini-3-pro"    = "Model B"
    "kimi-k2"         = "Model D"
}

# Determine the last used row in the sheet (data starts at row 2, header at row 1)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = $ws.UsedRange.Rows.Count }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value()
    if ($null -ne $current -and $mapping.ContainsKey($current)) {
        $cell.Value = $mapping[$current]
    }
}
